$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.029905206023092, 1.04245337638738, 1.029388595231997, 1.048110510162208, 1.034225811694147, 1.035049396683044, 1.045230020795166, 1.032202650128811, 1.050871269485186, 1.036519285548159)
    3 = @(1.02, 1.032772830575595, 1.043414495662642, 1.031899272155623, 1.049706511162192, 1.034452015615956, 1.037548471238826, 1.04600158099811, 1.034516885772671, 1.052277221066201, 1.039021909076458)
    4 = @(1.02, 1.034616769315617, 1.044032237793513, 1.033513183124531, 1.05073242178859, 1.034595437192268, 1.039154262702065, 1.046496058063569, 1.036003453037635, 1.053179691866135, 1.040629980947761)
    5 = @(1.02, 1.035389261183599, 1.044290954935468, 1.034189185804352, 1.051162117410158, 1.034655031764992, 1.03982670665717, 1.046702808737662, 1.036625861268074, 1.053557380616171, 1.041303379850413)
    6 = @(1.02, 1.035518809896693, 1.044334337570396, 1.034302545824282, 1.051234172477912, 1.034664997081711, 1.039939460667169, 1.046737457411951, 1.036730218978437, 1.053620696847821, 1.04141629398404)
    7 = @(1.02, 1.034627101894855, 1.044035698620732, 1.03352222558203, 1.050738169634745, 1.034596236240087, 1.03916325817696, 1.046498825085842, 1.036011779600371, 1.053184745238929, 1.040638989197263)
    8 = @(1.02, 1.030876784060089, 1.042779063413478, 1.030239339667297, 1.048651315937416, 1.034302872184683, 1.035896347512111, 1.045491771049145, 1.032987052528538, 1.051347942623944, 1.037367439144551)
    9 = @(1.02, 1.024175765070426, 1.040532079170092, 1.024369699759634, 1.044920364235416, 1.033763087697685, 1.030050114205575, 1.043679943256614, 1.027570689760299, 1.048054114613087, 1.031512903516147)
    10 = @(1.02, 1.019641129920965, 1.039011114839104, 1.020395202881001, 1.04239485135655, 1.033387495338198, 1.026087950884335, 1.042445989439899, 1.023897498012934, 1.045817756401392, 1.027545113468524)
    11 = @(1.02, 1.017660458292697, 1.038346834019006, 1.018658634070494, 1.041291724228018, 1.033221040019471, 1.024355928084862, 1.041905259218235, 1.02229123836904, 1.044839320702631, 1.025810630998009)
    12 = @(1.02, 1.016922070858163, 1.038099214286916, 1.018011165627941, 1.040880493968373, 1.033158629257308, 1.023710027205955, 1.041703422984996, 1.021692153126457, 1.044474330889916, 1.025163812865769)
    13 = @(1.02, 1.01708058031561, 1.038152369555341, 1.018150161184614, 1.040968771918835, 1.033172043034136, 1.023848692018978, 1.041746762447812, 1.021820771131214, 1.044552693494984, 1.025302674598731)
    14 = @(1.02, 1.017599478183953, 1.03832638369336, 1.018605164205062, 1.041257762245263, 1.033215893041446, 1.024302590273485, 1.041888595585226, 1.022241768245003, 1.044809182499372, 1.025757217440822)
    15 = @(1.02, 1.017918830352268, 1.038433482798276, 1.018885181983781, 1.041435621340426, 1.03324283314287, 1.024581911281445, 1.041975852525937, 1.022500831418701, 1.044967006506766, 1.026036935116652)
    16 = @(1.02, 1.019772210711609, 1.039055079248084, 1.020510117415657, 1.0424678571434, 1.033398461263366, 1.026202546835375, 1.042481738931503, 1.024003761394404, 1.045882476141583, 1.027659872158959)
    17 = @(1.02, 1.020930125987429, 1.03944345146282, 1.021525162492551, 1.04311276105841, 1.0334950542015, 1.027214680846042, 1.042797335046335, 1.024942236388731, 1.046453999241556, 1.028673443516014)
    18 = @(1.02, 1.021603871191, 1.039669434534346, 1.022115723401003, 1.043488002831946, 1.033551027057181, 1.027803467236154, 1.042980798884484, 1.025488119728216, 1.046786389467548, 1.029263066050335)
    19 = @(1.02, 1.021833324553149, 1.039746396618161, 1.022316838000931, 1.04361579585441, 1.033570050137153, 1.028003963837764, 1.043043251092187, 1.02567399766788, 1.046899562725991, 1.029463847380117)
    20 = @(1.02, 1.020806063741838, 1.039401839618011, 1.021416413305802, 1.043043664414613, 1.033484728830272, 1.027106251918156, 1.042763538642224, 1.024841704023308, 1.04639278079879, 1.028564860606613)
    21 = @(1.02, 1.017446750448149, 1.038275165226816, 1.018471244939126, 1.04117270299529, 1.033202996432195, 1.024168999688149, 1.041846856617615, 1.02211786337217, 1.044733696148333, 1.025623437141532)
    22 = @(1.02, 1.015319071355637, 1.037561699762283, 1.016605399262648, 1.039987771924169, 1.03302248971974, 1.022307430109942, 1.041264793079611, 1.020391065031226, 1.043681544494518, 1.023759223920999)
    23 = @(1.02, 1.016448502625923, 1.037940410023395, 1.017595885999345, 1.040616754553906, 1.033118501856372, 1.02329571703134, 1.041573904104397, 1.021307849122722, 1.044240178955292, 1.024748914323184)
    24 = @(1.02, 1.020862127214527, 1.03942064391155, 1.0214655570053, 1.043074889053695, 1.033489395558376, 1.027155251142865, 1.042778811691287, 1.024887134943165, 1.046420445771741, 1.028613929415842)
    25 = @(1.02, 1.025919616574523, 1.041116950317891, 1.02589763533318, 1.045891479363127, 1.033905379982178, 1.0315725689915, 1.044152863158772, 1.028981613622999, 1.048912628104409, 1.033037520362495)
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $rowVals[$i]
    }
}

Write-Host "Done writing $($data.Keys.Count) rows"
